$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 3770.2979
$ws.Range("I15").Value = 3770.2979
$ws.Range("K15").Value = 11310.8937
$ws.Range("M15").Value = -11141.8937
# Row 112
$ws.Range("H112").Value = 2024.3478
$ws.Range("I112").Value = 1539
$ws.Range("J112").Value = 2070.5715
$ws.Range("K112").Value = 4617
$ws.Range("L112").Value = 6211.7145
$ws.Range("M112").Value = -3509
$ws.Range("N112").Value = -8427.7145
# Row 132
$ws.Range("H132").Value = 5500.029
$ws.Range("I132").Value = 5609.5938
$ws.Range("K132").Value = 16828.7814
$ws.Range("M132").Value = -14298.7814
# Row 137
$ws.Range("H137").Value = 1730462.2
$ws.Range("I137").Value = 5001630.5
$ws.Range("K137").Value = 15004891.5
$ws.Range("M137").Value = -15002341.5
# Row 138
$ws.Range("H138").Value = 4058.7896
$ws.Range("I138").Value = 5296.636
$ws.Range("J138").Value = 3762.7827
$ws.Range("K138").Value = 15889.908
$ws.Range("L138").Value = 11288.3481
$ws.Range("M138").Value = -10749.908
$ws.Range("N138").Value = -21568.3481

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 989.7917
$ws.Range("I2").Value = 867.9
$ws.Range("J2").Value = 1599.25
$ws.Range("K2").Value = 867.9
$ws.Range("L2").Value = 1599.25
$ws.Range("M2").Value = -754.9
$ws.Range("N2").Value = -1825.25
# Row 34
$ws.Range("H34").Value = 23025
$ws.Range("I34").Value = 23025
$ws.Range("K34").Value = 23025
$ws.Range("M34").Value = -22754
# Row 61
$ws.Range("H61").Value = 2816.4075
$ws.Range("I61").Value = 2201.4
$ws.Range("J61").Value = 4573.5713
$ws.Range("K61").Value = 2201.4
$ws.Range("L61").Value = 4573.5713
$ws.Range("M61").Value = -1989.4
$ws.Range("N61").Value = -4997.5713
# Row 74
$ws.Range("H74").Value = 429796.7
$ws.Range("I74").Value = 795322.5600000001
$ws.Range("J74").Value = 3349.8333
$ws.Range("K74").Value = 795322.5600000001
$ws.Range("L74").Value = 3349.8333
$ws.Range("M74").Value = -794448.5600000001
$ws.Range("N74").Value = -5097.8333
# Row 77
$ws.Range("H77").Value = 429796.7
$ws.Range("I77").Value = 795322.5600000001
$ws.Range("J77").Value = 3349.8333
$ws.Range("K77").Value = 3976612.8
$ws.Range("L77").Value = 16749.1665
$ws.Range("M77").Value = -3972244.8
$ws.Range("N77").Value = -25485.1665
# Row 101
$ws.Range("H101").Value = 68999.5
$ws.Range("J101").Value = 68999.5
$ws.Range("L101").Value = 68999.5
$ws.Range("N101").Value = -75489.5
# Row 116
$ws.Range("H116").Value = 989.7917
$ws.Range("I116").Value = 867.9
$ws.Range("J116").Value = 1599.25
$ws.Range("K116").Value = 867.9
$ws.Range("L116").Value = 1599.25
$ws.Range("M116").Value = 1426.1
$ws.Range("N116").Value = -6187.25
# Row 128
$ws.Range("H128").Value = 100000
$ws.Range("J128").Value = 100000
$ws.Range("L128").Value = 100000
$ws.Range("N128").Value = -109960
# Row 132
$ws.Range("H132").Value = 2200.608
$ws.Range("I132").Value = 1659.8462
$ws.Range("J132").Value = 3958.0833
$ws.Range("K132").Value = 4979.5386
$ws.Range("L132").Value = 11874.2499
$ws.Range("M132").Value = -2449.5386
$ws.Range("N132").Value = -16934.2499
# Row 136
$ws.Range("H136").Value = 2816.4075
$ws.Range("I136").Value = 2201.4
$ws.Range("J136").Value = 4573.5713
$ws.Range("K136").Value = 6604.200000000001
$ws.Range("L136").Value = 13720.7139
$ws.Range("M136").Value = -4054.200000000001
$ws.Range("N136").Value = -18820.7139

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 989.7917
$ws.Range("I3").Value = 867.9
$ws.Range("J3").Value = 1599.25
$ws.Range("K3").Value = 867.9
$ws.Range("L3").Value = 1599.25
$ws.Range("M3").Value = -753.9
$ws.Range("N3").Value = -1827.25
# Row 75
$ws.Range("H75").Value = 20000
$ws.Range("I75").Value = 20000
$ws.Range("K75").Value = 20000
$ws.Range("M75").Value = -19064
# Row 78
$ws.Range("H78").Value = 20000
$ws.Range("I78").Value = 20000
$ws.Range("K78").Value = 60000
$ws.Range("M78").Value = -55320
# Row 80
$ws.Range("H80").Value = 933.2
$ws.Range("I80").Value = 903
$ws.Range("J80").Value = 940.75
$ws.Range("K80").Value = 903
$ws.Range("L80").Value = 940.75
$ws.Range("M80").Value = 95
$ws.Range("N80").Value = -2936.75
# Row 83
$ws.Range("H83").Value = 933.2
$ws.Range("I83").Value = 903
$ws.Range("J83").Value = 940.75
$ws.Range("K83").Value = 4515
$ws.Range("L83").Value = 4703.75
$ws.Range("M83").Value = 477
$ws.Range("N83").Value = -14687.75
# Row 86
$ws.Range("H86").Value = 1728.3334
$ws.Range("I86").Value = 1446
$ws.Range("J86").Value = 2187.125
$ws.Range("K86").Value = 1446
$ws.Range("L86").Value = 2187.125
$ws.Range("M86").Value = -323
$ws.Range("N86").Value = -4433.125
# Row 89
$ws.Range("H89").Value = 1728.3334
$ws.Range("I89").Value = 1446
$ws.Range("J89").Value = 2187.125
$ws.Range("K89").Value = 7230
$ws.Range("L89").Value = 10935.625
$ws.Range("M89").Value = -1614
$ws.Range("N89").Value = -22167.625
# Row 94
$ws.Range("H94").Value = 90910400
$ws.Range("I94").Value = 166666990
$ws.Range("J94").Value = 2493.8
$ws.Range("K94").Value = 166666990
$ws.Range("L94").Value = 2493.8
$ws.Range("M94").Value = -166666539
$ws.Range("N94").Value = -3395.8
# Row 105
$ws.Range("H105").Value = 9287518
$ws.Range("I105").Value = 557110.8
$ws.Range("J105").Value = 25002252
$ws.Range("K105").Value = 557110.8
$ws.Range("L105").Value = 25002252
$ws.Range("M105").Value = -555363.8
$ws.Range("N105").Value = -25005746
# Row 107
$ws.Range("H107").Value = 1673.1052
$ws.Range("I107").Value = 1449.5714
$ws.Range("K107").Value = 1449.5714
$ws.Range("M107").Value = 470.4286
# Row 134
$ws.Range("H134").Value = 2858.392
$ws.Range("I134").Value = 2662.3057
$ws.Range("K134").Value = 7986.9171
$ws.Range("M134").Value = -5451.9171

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4190.9688
$ws.Range("I31").Value = 3057.2856
$ws.Range("J31").Value = 6355.273
$ws.Range("K31").Value = 3057.2856
$ws.Range("L31").Value = 6355.273
$ws.Range("M31").Value = -2762.2856
$ws.Range("N31").Value = -6945.273
# Row 34
$ws.Range("H34").Value = 4190.9688
$ws.Range("I34").Value = 3057.2856
$ws.Range("J34").Value = 6355.273
$ws.Range("K34").Value = 3057.2856
$ws.Range("L34").Value = 6355.273
$ws.Range("M34").Value = -2855.2856
$ws.Range("N34").Value = -6759.273
# Row 58
$ws.Range("H58").Value = 4116.857
$ws.Range("I58").Value = 3940
$ws.Range("K58").Value = 3940
$ws.Range("M58").Value = -3737
# Row 122
$ws.Range("H122").Value = 2590
$ws.Range("I122").Value = 2255.9546
$ws.Range("K122").Value = 6767.8638
$ws.Range("M122").Value = -4317.8638
# Row 136
$ws.Range("H136").Value = 4116.857
$ws.Range("I136").Value = 3940
$ws.Range("K136").Value = 11820
$ws.Range("M136").Value = -9270

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 935.125
$ws.Range("I5").Value = 528.4545000000001
$ws.Range("K5").Value = 1585.3635
$ws.Range("M5").Value = -1473.3635
# Row 12
$ws.Range("H12").Value = 223.72728
$ws.Range("J12").Value = 282.85715
$ws.Range("L12").Value = 848.5714499999999
$ws.Range("N12").Value = -1194.57145
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
# Row 117
$ws.Range("H117").Value = 2115.0833
$ws.Range("J117").Value = 1909
$ws.Range("L117").Value = 5727
$ws.Range("N117").Value = -12611
# Row 121
$ws.Range("H121").Value = 150453.62
$ws.Range("I121").Value = 33533.332
$ws.Range("J121").Value = 220605.8
$ws.Range("K121").Value = 100599.996
$ws.Range("L121").Value = 661817.3999999999
$ws.Range("M121").Value = -99289.99600000001
$ws.Range("N121").Value = -664437.3999999999
# Row 125
$ws.Range("H125").Value = 5166.5
$ws.Range("I125").Value = 3333
$ws.Range("J125").Value = 7000
$ws.Range("K125").Value = 9999
$ws.Range("L125").Value = 21000
$ws.Range("M125").Value = -5079
$ws.Range("N125").Value = -30840
# Row 135
$ws.Range("H135").Value = 935.125
$ws.Range("I135").Value = 528.4545000000001
$ws.Range("K135").Value = 4756.0905
$ws.Range("M135").Value = -2221.0905

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 2090.158
$ws.Range("J113").Value = 2439.6667
$ws.Range("L113").Value = 2439.6667
$ws.Range("N113").Value = -6779.6667
# Row 132
$ws.Range("H132").Value = 2709.8125
$ws.Range("I132").Value = 1898.4546
$ws.Range("K132").Value = 5695.3638
$ws.Range("M132").Value = -3165.3638

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 696.41174
$ws.Range("I55").Value = 361
$ws.Range("K55").Value = 361
$ws.Range("M55").Value = -188
# Row 61
$ws.Range("H61").Value = 16919.562
$ws.Range("I61").Value = 3475.3333
$ws.Range("K61").Value = 3475.3333
$ws.Range("M61").Value = -3273.3333
# Row 100
$ws.Range("H100").Value = 1579.2
$ws.Range("I100").Value = 1366.3334
$ws.Range("K100").Value = 1366.3334
$ws.Range("M100").Value = -825.3334
# Row 113
$ws.Range("H113").Value = 16919.562
$ws.Range("I113").Value = 3475.3333
$ws.Range("K113").Value = 3475.3333
$ws.Range("M113").Value = -1305.3333
# Row 125
$ws.Range("H125").Value = 59996.5
$ws.Range("J125").Value = 59996.5
$ws.Range("L125").Value = 59996.5
$ws.Range("N125").Value = -69836.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
# Row 43
$ws.Range("H43").Value = 54030
$ws.Range("I43").Value = 54030
$ws.Range("K43").Value = 54030
$ws.Range("M43").Value = -53881
# Row 113
$ws.Range("H113").Value = 1523.2084
$ws.Range("J113").Value = 1565.4445
$ws.Range("L113").Value = 4696.333500000001
$ws.Range("N113").Value = -9036.333500000001
# Row 132
$ws.Range("H132").Value = 3889.55
$ws.Range("I132").Value = 3684.3076
$ws.Range("J132").Value = 4270.7144
$ws.Range("K132").Value = 11052.9228
$ws.Range("L132").Value = 12812.1432
$ws.Range("M132").Value = -8522.9228
